$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# The sheet held 2 periods' worth of debtor rows (rows 16-27, 4 workers,
# up to 4 periods each). The new data set is "parte 1" of a refreshed
# account-statement database: 11 workers, a single new period (2508),
# at the (higher) 56.940 / 1.423.500 salary bracket. That means one
# data row per worker instead of several, so the old 12-row block
# (16-27) shrinks to 11 rows (16-26): delete the old row 26, then
# overwrite rows 16-26 with the new workers. Row 27 (with the
# "closing" thicker-bottom-border style) naturally becomes row 26.
# ------------------------------------------------------------------

$ws.Rows.Item(26).Delete()

# New workers for period 2508 (all at 56.940 salario / 1.423.500 base)
$workers = @(
  @{ Tipo = "CC";  Doc = "1047411589"; Nombre = "MARCOS LEONARDO MORA MANJARREZ" },
  @{ Tipo = "CC";  Doc = "1043641747"; Nombre = "MAIRA ALEJANDRA AGUILAR FUENTES" },
  @{ Tipo = "CC";  Doc = "1041978941"; Nombre = "RAFAEL DAVID MONTERO OROZCO" },
  @{ Tipo = "CC";  Doc = "1128061759"; Nombre = "MANUEL ESTEBAN SIMANCAS AGAMEZ" },
  @{ Tipo = "PPT"; Doc = "6068342";    Nombre = "EDUARDO ANTONIO CHACIN MOLERO" },
  @{ Tipo = "CC";  Doc = "1002191088"; Nombre = "YON JAIRO POLO CARDONA" },
  @{ Tipo = "CC";  Doc = "1062961134"; Nombre = "LUIS FERNANDO ELLES ARGEL" },
  @{ Tipo = "CC";  Doc = "45505966";   Nombre = "EDITH SUSANA CASTILLO MUOZ" },
  @{ Tipo = "CC";  Doc = "64698775";   Nombre = "YESENIA YANETH AVILEZ NUEZ" },
  @{ Tipo = "CC";  Doc = "1048435905"; Nombre = "LUIS MIGUEL ALVAREZ VELASQUEZ" },
  @{ Tipo = "CC";  Doc = "1051443001"; Nombre = "JULIO DE JESUS ACEVEDO DIAZ" }
)

$row = 16
foreach ($w in $workers) {
  $ws.Cells.Item($row, 2).Value = $w.Tipo
  $ws.Cells.Item($row, 3).Value = $w.Doc
  $ws.Cells.Item($row, 4).Value = $w.Nombre
  $ws.Cells.Item($row, 5).Value = "2508"
  $ws.Cells.Item($row, 6).Value = 56940
  $ws.Cells.Item($row, 7).Value = 1423500
  $row = $row + 1
}

# Summary header: 11 workers, 1 period, total mora = 11 * 56940
$ws.Cells.Item(11, 5).Value = 626340
$ws.Cells.Item(13, 3).Value = 11
$ws.Cells.Item(13, 6).Value = 1
